$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column so numeric-looking strings
# (e.g. "1.001") are not auto-converted to numbers by the COM layer.
$priceCells = @(
"D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D50", "D51"
)
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.440.58"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "1.678.11"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "216.74"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.2700"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").Value = "0.06399"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").Value = "21.73"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("D11").Value = "0.07821"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "4.513"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "1.659.86"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "0.5564"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").Value = "0.0₅8320"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "65.61"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "26.480.37"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "4.735"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "193.52"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "6.341"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "142.22"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "0.1289"
$ws.Range("E25").Value = "  +5.74%  "
$ws.Range("D26").Value = "7.397"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "16.22"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "1.438"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").Value = "0.06243"
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("D30").Value = "1.275"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "3.609"
$ws.Range("E31").Value = "  +4.71%  "
$ws.Range("D32").Value = "3.448"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "1.678"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").Value = "1.008"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "0.6132"
$ws.Range("E35").Value = "  +6.84%  "
$ws.Range("D36").Value = "2.427"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "2.783"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01629"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "6.126"
$ws.Range("E39").Value = "  +7.46%  "
$ws.Range("D40").Value = "1.084.77"
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("D41").Value = "0.8650"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D43").Value = "100.24"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "1.822.38"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("D45").Value = "57.15"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "8.143"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "1.475"
$ws.Range("E50").Value = "  +6.11%  "
$ws.Range("D51").Value = "6.021"
$ws.Range("E51").Value = "  +1.81%  "

foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}
